# VIN-379 - elaborando a jornada do usuario
# Updates the "Referencias" section:
#  - Cleans up the "Ozbun" reference (removes spell-check proofErr run
#    splitting) and adds a trailing period.
#  - Adds two new references (Pressman et al.; Sommerville) into the
#    blank paragraphs that follow it.

$d = $word.ActiveDocument

function New-CleanParagraphWithRuns {
    param(
        $AfterParagraph,
        [string[]] $RunTexts
    )

    # Insert a brand-new, plain paragraph (same style as its neighbours)
    # right after $AfterParagraph.
    $AfterParagraph.Range.InsertParagraphAfter() | Out-Null
    $newIndex = $AfterParagraph.Index + 1
    $newPara = $d.Paragraphs.Item($newIndex)

    # Write the first chunk of text without tracked changes so the
    # paragraph's implicit run is simply filled in (no stray empty run).
    $r = $newPara.Range
    $r.MoveEnd(1, -1)
    $r.Text = $RunTexts[0]

    # Any remaining chunks are appended with Track Changes enabled so
    # that each one lands in its own <w:r> (otherwise runs with
    # identical formatting are coalesced on save). Accepting the
    # revisions afterwards removes the ins/del markup while keeping the
    # run boundaries intact.
    if ($RunTexts.Count -gt 1) {
        $d.TrackRevisions = $true
        for ($i = 1; $i -lt $RunTexts.Count; $i++) {
            $curPara = $d.Paragraphs.Item($newIndex)
            $endPos = $curPara.Range.End - 1
            $tail = $d.Range($endPos, $endPos)
            $tail.InsertAfter($RunTexts[$i])
        }
        $d.TrackRevisions = $false
        $d.AcceptAllRevisions()
    }

    return $d.Paragraphs.Item($newIndex)
}

# --- 1. Fix the "Ozbun" reference paragraph -------------------------------
$beforeOzbun = $d.Paragraphs.Item(12)
$ozbunPara = $d.Paragraphs.Item(13)
$delRange = $d.Range($ozbunPara.Range.Start, $ozbunPara.Range.End)
$delRange.Delete()

$ozbunText1 = "Ozbun, T. Wine industry in Brazil " + [char]0x2013 + " statistics & facts. Statista, Jan 10, 2024. Disponível em: https://www.statista.com/topics/5228/wine-industry-in-brazil/#topicOverview. Acesso em: 22 mar. 2024"
$ozbunText2 = "."
New-CleanParagraphWithRuns $beforeOzbun @($ozbunText1, $ozbunText2) | Out-Null

# --- 2. Add the "Pressman" reference in the first blank paragraph --------
$ozbunParaNow = $d.Paragraphs.Item(13)
$blank1 = $d.Paragraphs.Item(14)

$pressman1 = "Pressman, R. S., Maxim, B, R. Engenharia de software: uma abordagem Profissional. "
$pressman2 = "9. ed. " + [char]0x2013 + " Porto Alegre: AMGH, 2021"
$pressman3 = "."
New-CleanParagraphWithRuns $blank1 @($pressman1, $pressman2, $pressman3) | Out-Null

# --- 3. Add the "Sommerville" reference in the next blank paragraph ------
$blank2 = $d.Paragraphs.Item(16)
$sommerville1 = "Sommerville, I. Engenharia de software. 10 ed, São Paulo. Pearson Education do Brasil, 2018."
New-CleanParagraphWithRuns $blank2 @($sommerville1) | Out-Null
